# Fix header labels on the existing sheets so they describe the PO
# quantity metric explicitly instead of the generic "Requested quantity".
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$ws2 = $wb.Worksheets.Item(2)   # "Monthly Trend"

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add a new "PO Forecast" sheet after the existing sheets.
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "PO Forecast"

# Reuse the existing header/date cell formatting so the new sheet matches
# the look of the other two (bold+bordered header row, date number format
# on column A).
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws3.Range("A2:A12").PasteSpecial(-4122)

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

$ws3.Range("A2").Value = 45578.99999999999
$ws3.Range("B2").Value = 81
$ws3.Range("C2").Value = -15.16915423984335
$ws3.Range("D2").Value = 173.6189205461679
$ws3.Range("A3").Value = 45592.99999999999
$ws3.Range("B3").Value = 117
$ws3.Range("C3").Value = 25.27364854429067
$ws3.Range("D3").Value = 203.7519592282909
$ws3.Range("A4").Value = 45613.99999999999
$ws3.Range("B4").Value = 171
$ws3.Range("C4").Value = 76.60921765316468
$ws3.Range("D4").Value = 259.3421716364832
$ws3.Range("A5").Value = 45620.99999999999
$ws3.Range("B5").Value = 189
$ws3.Range("C5").Value = 95.02748697043964
$ws3.Range("D5").Value = 276.2834399995787
$ws3.Range("A6").Value = 45627.99999999999
$ws3.Range("B6").Value = 206
$ws3.Range("C6").Value = 107.7602535534828
$ws3.Range("D6").Value = 304.0176548993385
$ws3.Range("A7").Value = 45634.99999999999
$ws3.Range("B7").Value = 224
$ws3.Range("C7").Value = 125.6455961322709
$ws3.Range("D7").Value = 314.4003176204922
$ws3.Range("A8").Value = 45641.99999999999
$ws3.Range("B8").Value = 242
$ws3.Range("C8").Value = 151.1644262030852
$ws3.Range("D8").Value = 330.5202706582909
$ws3.Range("A9").Value = 45648.99999999999
$ws3.Range("B9").Value = 260
$ws3.Range("C9").Value = 165.8831788082284
$ws3.Range("D9").Value = 351.5150799438555
$ws3.Range("A10").Value = 45655.99999999999
$ws3.Range("B10").Value = 278
$ws3.Range("C10").Value = 179.1352999087144
$ws3.Range("D10").Value = 378.3258162091458
$ws3.Range("A11").Value = 45662.99999999999
$ws3.Range("B11").Value = 296
$ws3.Range("C11").Value = 204.0176197342864
$ws3.Range("D11").Value = 387.3470848302016
$ws3.Range("A12").Value = 45669.99999999999
$ws3.Range("B12").Value = 314
$ws3.Range("C12").Value = 222.0418464452689
$ws3.Range("D12").Value = 405.551502997226

# Keep the originally-active sheet selected (adding a sheet makes it active
# by default).
$ws1.Activate()

Write-Output "PO Forecast sheet populated"
